$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new bug-tracking row (row 9) with bug #5: the respawn-at-checkpoint bug.
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = 42480
$ws.Range("C9").Value = "Charles"
$ws.Range("D9").Value = "Code - Function"
$ws.Range("E9").Value = "High"
$ws.Range("F9").Value = "Charles"
$ws.Range("G9").Value = "Player doesn't repsawn at designated checkpoint after death"

# This bug hasn't been fixed yet (no Date Fixed), so flag the row with Excel's
# built-in "Bad" cell style (red), matching how rows 6-8 use "Check Cell" (gray,
# fixed) and row 5 uses "Neutral" (yellow, in progress).
$row9 = $ws.Range("A9:L9")
$row9.Style = "Bad"
$row9.HorizontalAlignment = -4108
$row9.VerticalAlignment = -4108
$row9.Borders.LineStyle = 1
$row9.Borders.Weight = 2

$ws.Range("G9:K9").WrapText = $true
$ws.Range("B9").NumberFormat = "m/d/yyyy"

# Merge the summary cells for this row, matching the other bug rows.
$ws.Range("G9:K9").Merge()

# Update the view: scroll so row 4 is at top, and move the selection.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("N8").Select()
